# Add 2022-Q3 data
#
# 1) "总计" (summary) sheet: insert a new row for 2022-Q3 right after the
#    header row, pushing the existing quarters down by one row.
# 2) Insert a brand-new "2022-Q3" worksheet (built from a copy of the
#    "2021-Q2" sheet so that styles/page-setup match the rest of the
#    workbook) positioned right after "总计", containing the two fund rows.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, [string]$text) {
    # Force a cell to hold literal TEXT even when the text looks numeric
    # (e.g. "010690" or "0.55"), without disturbing the cell's existing
    # number format / border / font (same-cell self-paste keeps style put).
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

# ---------------------------------------------------------------------
# 1) "总计" sheet - shift rows 2..4 down to 3..5, then write the new
#    2022-Q3 row into row 2.
# ---------------------------------------------------------------------
$sheetTotal = $wb.Worksheets.Item("总计")

$existing = @()
for ($r = 2; $r -le 4; $r++) {
    $rowVals = @(
        $sheetTotal.Cells.Item($r, 1).Value2,
        $sheetTotal.Cells.Item($r, 2).Value2,
        $sheetTotal.Cells.Item($r, 3).Value2,
        $sheetTotal.Cells.Item($r, 4).Value2
    )
    $existing += ,$rowVals
}

for ($i = 2; $i -ge 0; $i--) {
    $destRow = $i + 3
    $sheetTotal.Cells.Item($destRow, 1).Value = $i + 1
    $sheetTotal.Cells.Item($destRow, 2).Value = $existing[$i][1]
    $sheetTotal.Cells.Item($destRow, 3).Value = $existing[$i][2]
    $sheetTotal.Cells.Item($destRow, 4).Value = $existing[$i][3]
}

# row 5 is brand new -> give A5 the same style as A4 (bold/bordered index)
$sheetTotal.Cells.Item(4, 1).Copy()
$sheetTotal.Cells.Item(5, 1).PasteSpecial(-4122)  # xlPasteFormats

$sheetTotal.Cells.Item(2, 1).Value = 0
$sheetTotal.Cells.Item(2, 2).Value = "2022-Q3"
$sheetTotal.Cells.Item(2, 3).Value = 2
$sheetTotal.Cells.Item(2, 4).Value = 0.05

# ---------------------------------------------------------------------
# 2) New "2022-Q3" worksheet, inserted right after "总计".
#    Built as a copy of "2021-Q2" so sheetPr/pageMargins/column styles
#    match the other quarter sheets exactly, then trimmed/overwritten.
# ---------------------------------------------------------------------
$sheetQ2 = $wb.Worksheets.Item("2021-Q2")
$sheetQ2.Copy($null, $sheetTotal)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# Only 2 data rows are needed (the template sheet had 5)
$newSheet.Rows("4:6").Delete()

# Header row
$newSheet.Range("B1").Value2 = "基金代码"
$newSheet.Range("C1").Value2 = "基金名称"
$newSheet.Range("D1").Value2 = "基金规模"
$newSheet.Range("E1").Value2 = "股票总仓位"
$newSheet.Range("F1").Value2 = "仓位占比"
$newSheet.Range("G1").Value2 = "持有市值(亿元)"
$newSheet.Range("H1").Value2 = "仓位排名"

# Row 2 - 010690
$newSheet.Range("A2").Value2 = 0
Set-TextValue $newSheet.Range("B2") "010690"
Set-TextValue $newSheet.Range("C2") "万家互联互通核心资产量化策略混合A"
Set-TextValue $newSheet.Range("D2") "0.55"
Set-TextValue $newSheet.Range("E2") "92.41"
Set-TextValue $newSheet.Range("F2") "6.98"
Set-TextValue $newSheet.Range("G2") "0.0384"
$newSheet.Range("H2").Value2 = 8

# Row 3 - 010691
$newSheet.Range("A3").Value2 = 1
Set-TextValue $newSheet.Range("B3") "010691"
Set-TextValue $newSheet.Range("C3") "万家互联互通核心资产量化策略混合C"
Set-TextValue $newSheet.Range("D3") "0.15"
Set-TextValue $newSheet.Range("E3") "92.41"
Set-TextValue $newSheet.Range("F3") "6.98"
Set-TextValue $newSheet.Range("G3") "0.0105"
$newSheet.Range("H3").Value2 = 8
